$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# --- Row 104: The Big Short ---
$ws.Range("A104").Value = "The Big Short"
$ws.Range("B104").Value = "Michael Lewis"

$ws.Range("C103").Copy()
$ws.Range("C104").PasteSpecial(-4122)
$ws.Range("C104").Value = 44018

$ws.Range("D103").Copy()
$ws.Range("D104").PasteSpecial(-4122)
$ws.Range("D104").Value = 44019

$ws.Range("E104").Value = "bonds;mortgage backed securities;shorting;wall street;finance;financial crisis"
$ws.Range("F104").Value = "Audio"
$ws.Range("G104").Value = "9 Hours 30 Mins"

# --- Row 105: What the Most Successful People Do At Work ---
$ws.Range("A105").Value = "What the Most Successful People Do At Work"
$ws.Range("B105").Value = "Laura Vanderkam"

$ws.Range("C103").Copy()
$ws.Range("C105").PasteSpecial(-4122)
$ws.Range("C105").Value = 44019

$ws.Range("D103").Copy()
$ws.Range("D105").PasteSpecial(-4122)
$ws.Range("D105").Value = 44019

$ws.Range("E105").Value = "success;time management;planning;career;self-improvement"
$ws.Range("F105").Value = "Audio"
$ws.Range("G105").Value = "2 Hours 4 Mins"

$ws.Range("A106").Select() | Out-Null
